# Update the "Forecast Comparison" sheet:
#  - insert a new "Week_Start_Date" column between "Week" and "ASIN"
#  - change the "Week" labels from W01.."W16" to W1.."W16"
#  - fill in weekly start dates (as text, not Excel dates)
#  - store is_holiday_week as boolean FALSE instead of numeric 0
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Make room for the new column (shifts ASIN..is_holiday_week right by one).
$ws.Columns("B:B").Insert()

$ws.Range("B1").Value = "Week_Start_Date"

$weekStarts = @(
    "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
    "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
    "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
    "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)

for ($i = 0; $i -lt 16; $i++) {
    $r = $i + 2

    # "W01" -> "W1", etc.
    $ws.Cells.Item($r, 1).Value = "W" + ($i + 1)

    # Write the week-start date as plain text (not an auto-converted date).
    $cell = $ws.Cells.Item($r, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $weekStarts[$i]
    $cell.ClearFormats()

    # is_holiday_week (column J after the insert) becomes a real boolean.
    $ws.Cells.Item($r, 10).Value = $false
}

# Update the dependent summary metric on the "Summary" sheet (stored as
# text, like the rest of column B on this sheet).
$summary = $wb.Worksheets.Item("Summary")
$summaryCell = $summary.Range("B10")
$summaryCell.NumberFormat = "@"
$summaryCell.Value = "257"
$summaryCell.ClearFormats()
